$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.066.23'
$ws.Range('E2').Value = '  -1.40%  '
$ws.Range('D3').Value = '2.689.63'
$ws.Range('E3').Value = '  -1.81%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '556.80'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.29'
$ws.Range('D6').ClearFormats()
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.80%  '
$ws.Range('E9').Value = '  -2.55%  '
$ws.Range('E10').Value = '  -1.98%  '
$ws.Range('E11').Value = '  -3.05%  '
$ws.Range('E12').Value = '  -5.30%  '
$ws.Range('D13').Value = '3.166.55'
$ws.Range('E13').Value = '  -1.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.63'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('D15').Value = '62.984.19'
$ws.Range('E15').Value = '  -1.27%  '
$ws.Range('D17').Value = '2.691.41'
$ws.Range('E17').Value = '  -1.88%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.01'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.54%  '
$ws.Range('E19').Value = '  -3.12%  '
$ws.Range('E20').Value = '  -2.03%  '
$ws.Range('E21').Value = '  -4.18%  '
$ws.Range('E22').Value = '  -0.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.58'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.14%  '
$ws.Range('E25').Value = '  -0.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.28'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.46'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +9.56%  '
$ws.Range('E29').Value = '  -4.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.32'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.16%  '
$ws.Range('E31').Value = '  -0.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '164.48'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.96'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.12%  '
$ws.Range('E34').Value = '  +1.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '19.59'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.56%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.81'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '362.12'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +4.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.54'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +3.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.967'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.04'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '38.58'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.22'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.48'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.99%  '
$ws.Range('E45').Value = '  -2.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.622'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.999'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('E49').Value = '  -2.31%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0976'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.62%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '129.90'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.37%  '
